$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 159 (pushes existing rows 159:166 down to 160:167)
$ws.Rows(159).Insert()

# Populate the newly inserted row 159 with the new weekly price record
$ws.Range("A159").Value = 10
$ws.Range("B159").Value = "Vega Modelo de Temuco"
$ws.Range("C159").Value = "La Araucanía"
$ws.Range("D159").Value = 44706
$ws.Range("E159").Value = 9
$ws.Range("F159").Value = 100114007
$ws.Range("G159").Value = "Jengibre"
$ws.Range("H159").Value = "Sin especificar"
$ws.Range("I159").Value = "Primera"
$ws.Range("J159").Value = 40
$ws.Range("K159").Value = 20000
$ws.Range("L159").Value = 20000
$ws.Range("M159").Value = 20000
$ws.Range("N159").Value = "$/caja 13 kilos"
$ws.Range("O159").Value = "Perú"
$ws.Range("P159").Value = 1538
$ws.Range("Q159").Value = 13
$ws.Range("R159").Value = "Hortaliza"
